# Applies the attendance_reports sync update to
# Y3_B2526_Urogenital_session_analysis.xlsx (sheet: "Session Analysis Results")
#
# Summary of changes:
#  1. Reorder the recorded-by email lists in several rows (same people, new order).
#  2. Update the "Class Statistics" summary numbers (L6:L10).
#  3. Update the per-group breakdown numbers (O16:S16, O18:S19).
#  4. Three sessions that were previously "Not Recorded"/"Pending" now have
#     attendance recorded (rows 24, 66, 85): fill in Recorded By / Students /
#     Status, and restyle those rows to the normal "Recorded" look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($addr, $text) {
    $ws.Range($addr).Value2 = $text
}

# Excel's COM layer auto-parses strings that look like "31.6%" into the
# numeric value 0.316 with a Percent number format. The source workbook
# stores these as literal text, so force a Text number format first to
# keep them as plain strings.
function Set-TextPercent($addr, $text) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value2 = $text
}

# ---------------------------------------------------------------------------
# 1. Recorded-by email list reorderings (membership unchanged, order changed)
# ---------------------------------------------------------------------------

$group1 = "Veronia.rafat@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, heba@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
foreach ($addr in @("G2","G21","G40")) { Set-Text $addr $group1 }

Set-Text "G5" "AbeerRagheb@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg"

$group3 = "ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, eman.samir@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
foreach ($addr in @("G18","G37","G56","G75","G94","G113")) { Set-Text $addr $group3 }

$group4 = "yasmin.m.senosy@med.asu.edu.eg, eman.samir@med.asu.edu.eg, marinasorial@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
foreach ($addr in @("G19","G76","G95")) { Set-Text $addr $group4 }

$group5 = "yasmin.m.senosy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
foreach ($addr in @("G20","G38","G39","G57","G58","G77","G96","G115")) { Set-Text $addr $group5 }

Set-Text "G28" "yassmina.fattoh@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, dina.adel@med.asu.edu.eg"

$group7 = "NadaMohamed@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Fatmaelhady@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg"
foreach ($addr in @("G43","G100")) { Set-Text $addr $group7 }

$group8 = "Mohammedeltanany@med.asu.edu.eg, heba@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg"
foreach ($addr in @("G59","G78","G97")) { Set-Text $addr $group8 }

# ---------------------------------------------------------------------------
# 2. Class Statistics summary box (K5:L10)
# ---------------------------------------------------------------------------

$ws.Range("L6").Value2 = 36
$ws.Range("L7").Value2 = 2
$ws.Range("L8").Value2 = 76
Set-TextPercent "L9" "31.6%"
Set-TextPercent "L10" "42.7%"

# ---------------------------------------------------------------------------
# 3. Per-group breakdown table (rows 16-19)
# ---------------------------------------------------------------------------

# Row 16 - Year 3 / A2
$ws.Range("O16").Value2 = 7
$ws.Range("P16").Value2 = 0
Set-TextPercent "R16" "36.8%"
Set-TextPercent "S16" "44.9%"

# Row 18 - Year 3 / B1
$ws.Range("O18").Value2 = 6
$ws.Range("Q18").Value2 = 12
Set-TextPercent "R18" "31.6%"
Set-TextPercent "S18" "23.0%"

# Row 19 - Year 3 / B2
$ws.Range("O19").Value2 = 6
$ws.Range("P19").Value2 = 0
Set-TextPercent "R19" "31.6%"
Set-TextPercent "S19" "36.6%"

# ---------------------------------------------------------------------------
# 4. Sessions that moved from Not Recorded/Pending -> Recorded
#    (rows 24, 66, 85). Copy the usual "Recorded" row look (format only)
#    from row 2, then fill in the new Recorded By / Students / Status data.
# ---------------------------------------------------------------------------

$ws.Range("A2:I2").Copy()
$ws.Range("A24:I24").PasteSpecial(-4122)
$ws.Range("A66:I66").PasteSpecial(-4122)
$ws.Range("A85:I85").PasteSpecial(-4122)
$excel.CutCopyMode = 0

Set-Text "G24" "Fatmaelhady@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg"
Set-Text "H24" "15/204"
Set-Text "I24" "Recorded"

Set-Text "G66" "amira.m.ibrahim@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"
Set-Text "H66" "2/149"
Set-Text "I66" "Recorded"

Set-Text "G85" "amira.m.ibrahim@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"
Set-Text "H85" "97/227"
Set-Text "I85" "Recorded"
